$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to be forced to Text format so Excel does not re-interpret
# numeric-looking strings (e.g. "12.30", "0.999") and strip formatting/precision.
$textCells = @("D4", "D5", "D6", "D8", "D12", "D14", "D18", "D20", "D21", "D26", "D27", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D45", "D46", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values captured from the source diff.
$ws.Range('D2').Value = '65.816.36'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '2.954.06'
$ws.Range('E3').Value = '  -2.02%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '572.12'
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('D6').Value = '162.58'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.516'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').Value = '2.950.56'
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('E10').Value = '  -3.43%  '
$ws.Range('E11').Value = '  -4.60%  '
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('E13').Value = '  -2.76%  '
$ws.Range('D14').Value = '34.81'
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('E15').Value = '  -1.64%  '
$ws.Range('D16').Value = '65.623.08'
$ws.Range('E16').Value = '  -0.64%  '
$ws.Range('D17').Value = '3.441.42'
$ws.Range('D18').Value = '7.03'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D19').Value = '2.951.67'
$ws.Range('E19').Value = '  -1.96%  '
$ws.Range('D20').Value = '15.94'
$ws.Range('E20').Value = '  +13.95%  '
$ws.Range('D21').Value = '447.32'
$ws.Range('E21').Value = '  -2.43%  '
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('E23').Value = '  -1.16%  '
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('D26').Value = '12.30'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').Value = '10.05'
$ws.Range('E27').Value = '  -5.73%  '
$ws.Range('D29').Value = '2.53'
$ws.Range('E29').Value = '  +8.12%  '
$ws.Range('D30').Value = '8.12'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('D31').Value = '0.0000104'
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('D32').Value = '2.60'
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('E33').Value = '  +2.65%  '
$ws.Range('D34').Value = '27.29'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = '0.975'
$ws.Range('E36').Value = '  -2.12%  '
$ws.Range('D37').Value = '5.74'
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('D38').Value = '45.59'
$ws.Range('E38').Value = '  +4.01%  '
$ws.Range('D39').Value = '49.23'
$ws.Range('E39').Value = '  -1.44%  '
$ws.Range('D40').Value = '2.00'
$ws.Range('E40').Value = '  -6.86%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = '0.303'
$ws.Range('E41').Value = '  -1.42%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.122'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('E43').Value = '  -6.58%  '
$ws.Range('E44').Value = '  +0.69%  '
$ws.Range('D45').Value = '387.44'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '0.0353'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('D47').Value = '2.685.01'
$ws.Range('E47').Value = '  -4.08%  '
$ws.Range('D48').Value = '133.22'
$ws.Range('E48').Value = '  -1.57%  '
$ws.Range('D50').Value = '23.89'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').Value = '2.18'
$ws.Range('E51').Value = '  +1.11%  '
